$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set E24:E29 to 64 for existing rows 24-29
foreach ($r in 24..29) {
    $ws.Cells.Item($r, 5).Value = 64
}

# Add new rows 30-48
# Row 30
$ws.Cells.Item(30, 3).Value = 216066
$ws.Cells.Item(30, 5).Value = 64
$ws.Cells.Item(30, 6).Value = 3.0548000000000002
$ws.Cells.Item(30, 7).Value = 4

# Row 31
$ws.Cells.Item(31, 3).Value = 295087
$ws.Cells.Item(31, 5).Value = 64
$ws.Cells.Item(31, 6).Value = 2.0613999999999999
$ws.Cells.Item(31, 7).Value = 7

# Row 32
$ws.Cells.Item(32, 3).Value = 176035
$ws.Cells.Item(32, 5).Value = 64
$ws.Cells.Item(32, 6).Value = 1.6104000000000001
$ws.Cells.Item(32, 7).Value = 7

# Row 33
$ws.Cells.Item(33, 3).Value = 296059
$ws.Cells.Item(33, 5).Value = 64
$ws.Cells.Item(33, 6).Value = 2.7111000000000001
$ws.Cells.Item(33, 7).Value = 5

# Row 34
$ws.Cells.Item(34, 3).Value = 124084
$ws.Cells.Item(34, 5).Value = 64
$ws.Cells.Item(34, 6).Value = 1.0268999999999999
$ws.Cells.Item(34, 7).Value = 10

# Row 35
$ws.Cells.Item(35, 3).Value = 163014
$ws.Cells.Item(35, 5).Value = 64
$ws.Cells.Item(35, 6).Value = 2.3742000000000001
$ws.Cells.Item(35, 7).Value = 5

# Row 36
$ws.Cells.Item(36, 3).Value = 35070
$ws.Cells.Item(36, 5).Value = 64
$ws.Cells.Item(36, 6).Value = 2.3170000000000002
$ws.Cells.Item(36, 7).Value = 6

# Row 37
$ws.Cells.Item(37, 3).Value = 42049
$ws.Cells.Item(37, 5).Value = 64
$ws.Cells.Item(37, 6).Value = 2.2507999999999999
$ws.Cells.Item(37, 7).Value = 5

# Row 38
$ws.Cells.Item(38, 3).Value = 253036
$ws.Cells.Item(38, 5).Value = 64
$ws.Cells.Item(38, 6).Value = 2.5701999999999998
$ws.Cells.Item(38, 7).Value = 5

# Row 39
$ws.Cells.Item(39, 3).Value = 'test'
$ws.Cells.Item(39, 5).Value = 64
$ws.Cells.Item(39, 6).Value = 1.7512000000000001
$ws.Cells.Item(39, 7).Value = 7

# Row 40
$ws.Cells.Item(40, 3).Value = 188091
$ws.Cells.Item(40, 5).Value = 64
$ws.Cells.Item(40, 6).Value = 2.6861000000000002
$ws.Cells.Item(40, 7).Value = 4

# Row 41
$ws.Cells.Item(41, 3).Value = 45096
$ws.Cells.Item(41, 5).Value = 64
$ws.Cells.Item(41, 6).Value = 1.0215000000000001
$ws.Cells.Item(41, 7).Value = 10

# Row 42
$ws.Cells.Item(42, 3).Value = 60079
$ws.Cells.Item(42, 5).Value = 64
$ws.Cells.Item(42, 6).Value = 4.1104000000000003
$ws.Cells.Item(42, 7).Value = 4

# Row 43
$ws.Cells.Item(43, 3).Value = 'Plain'
$ws.Cells.Item(43, 5).Value = 64
$ws.Cells.Item(43, 6).Value = 0.82157999999999998
$ws.Cells.Item(43, 7).Value = 15

# Row 44
$ws.Cells.Item(44, 3).Value = 'showimage'
$ws.Cells.Item(44, 5).Value = 64
$ws.Cells.Item(44, 6).Value = 1.2588999999999999
$ws.Cells.Item(44, 7).Value = 7

# Row 45
$ws.Cells.Item(45, 3).Value = 'leucolinf'
$ws.Cells.Item(45, 5).Value = 64
$ws.Cells.Item(45, 6).Value = 4.5297999999999998
$ws.Cells.Item(45, 7).Value = 4

# Row 46
$ws.Cells.Item(46, 3).Value = 135069
$ws.Cells.Item(46, 5).Value = 64
$ws.Cells.Item(46, 6).Value = 2.6396000000000002
$ws.Cells.Item(46, 7).Value = 4

# Row 47
$ws.Cells.Item(47, 3).Value = 'linfocitos1'
$ws.Cells.Item(47, 5).Value = 64
$ws.Cells.Item(47, 6).Value = 5.4508999999999999
$ws.Cells.Item(47, 7).Value = 4

# Row 48
$ws.Cells.Item(48, 3).Value = '10472953_992127864183797_1990666493_n'
$ws.Cells.Item(48, 5).Value = 64
$ws.Cells.Item(48, 6).Value = 1.2169000000000001
$ws.Cells.Item(48, 7).Value = 8

# Update selection / view
$ws.Range("F51").Select()
